$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into the Price column (D) as literal text, even when
# it looks like a number Excel would otherwise "clean up" (stripping
# trailing zeros, e.g. "0.07700" -> 0.077). Cells that are not at risk are
# written directly so we don't touch their style (keeps style index 0).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row -> @(D_new, E_new) for rows whose Price and/or Volume(1h) changed.
# (rows 12, 18, 20, 27 only change column E; all others change D and E)
$updates = @{
    2  = @("28.058.55", "  +1.58%  ")
    3  = @("1.808.92",  "  +2.67%  ")
    4  = @("1.003",     "  -1.05%  ")
    5  = @("329.87",    "  +1.11%  ")
    6  = @("1.001",     "  -0.84%  ")
    7  = @("0.4441",    "  +3.79%  ")
    8  = @("0.3726",    "  +2.76%  ")
    9  = @("44.78",     "  -0.87%  ")
    10 = @("0.07700",   "  +4.26%  ")
    11 = @("1.119",     "  -0.32%  ")
    12 = @($null,       "  -0.95%  ")
    13 = @("21.97",     "  +0.91%  ")
    14 = @("6.297",     "  +2.50%  ")
    15 = @("7.463",     "  +2.47%  ")
    16 = @("1.813.79",  "  +2.70%  ")
    17 = @("93.70",     "  +12.79%  ")
    18 = @($null,       "  +1.61%  ")
    19 = @("0.06488",   "  +4.38%  ")
    20 = @($null,       "  -0.67%  ")
    21 = @("17.48",     "  +3.12%  ")
    22 = @("6.260",     "  +2.34%  ")
    23 = @("0.5343",    "  -2.37%  ")
    24 = @("28.110.67", "  +1.49%  ")
    25 = @("11.70",     "  +3.59%  ")
    26 = @("2.151",     "  -10.51%  ")
    27 = @($null,       "  +2.36%  ")
    28 = @("155.46",    "  +2.44%  ")
    29 = @("2.019.99",  "  +2.56%  ")
    30 = @("2.325",     "  -2.30%  ")
    31 = @("127.34",    "  -0.29%  ")
    32 = @("1.201",     "  -7.04%  ")
    33 = @("5.848",     "  +3.77%  ")
    34 = @("0.09228",   "  +1.91%  ")
    35 = @("3.669",     "  -6.54%  ")
    36 = @("13.08",     "  +5.22%  ")
    37 = @("0.02342",   "  +2.92%  ")
    38 = @("0.2171",    "  -0.02%  ")
    39 = @("5.170",     "  +2.06%  ")
    42 = @("1.193",     "  +0.64%  ")
    43 = @("8.072",     "  +1.56%  ")
    44 = @("1.001",     "  -0.78%  ")
    45 = @("13.99",     "  +0.80%  ")
    46 = @("1.390",     "  -2.60%  ")
    47 = @("0.6068",    "  +2.18%  ")
    48 = @("3.764",     "  -0.50%  ")
    49 = @("126.61",    "  +1.24%  ")
    50 = @("2.026",     "  +3.25%  ")
}

# Rows whose new Price text would otherwise be mangled by Excel's automatic
# number coercion (trailing zero stripped) if assigned the normal way.
$riskyRows = @(10, 17, 22, 25, 39, 46)

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals[0]) {
        $dCell = $ws.Cells.Item($row, 4)
        if ($riskyRows -contains $row) {
            Set-TextValue $dCell $vals[0]
        } else {
            $dCell.Value = $vals[0]
        }
    }
    $ws.Cells.Item($row, 5).Value = $vals[1]
}

# Rows 40 & 41 swap their Hedera / TheSandbox entries (with updated figures).
$ws.Cells.Item(40, 2).Value = "Hedera"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(40, 4).Value = "0.06201"
$ws.Cells.Item(40, 5).Value = "  +0.46%  "

$ws.Cells.Item(41, 2).Value = "TheSandbox"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(41, 4).Value = "0.6556"
$ws.Cells.Item(41, 5).Value = "  +0.60%  "

# Row 51: EOS dropped out of the top-50, replaced by Cronos.
$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(51, 4).Value = "0.06983"
$ws.Cells.Item(51, 5).Value = "  +1.30%  "
